# Update the "想去人数" (interested count) figures in the "展览" and
# "全部类型" worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 11163
    $ws.Range("F5").Value = 10403
    $ws.Range("F13").Value = 10445
    $ws.Range("F18").Value = 14
    $ws.Range("F21").Value = 11063
}
